# Apply updated crypto price/volume data to Sheet1 (columns D and E, rows 2-51)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.157.81'
$ws.Range("E2").Value = '  -0.15%  '
$ws.Range("D3").Value = '1.669.94'
$ws.Range("E3").Value = '  -0.63%  '
$ws.Range("E4").Value = '  -0.28%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '210.91'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.52%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5203'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.10%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.004'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.26%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2626'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.38%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06325'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.65%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.19'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.06%  '
$ws.Range("E11").Value = '  -0.88%  '
$ws.Range("D12").Value = '1.670.28'
$ws.Range("E12").Value = '  -0.78%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.443'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.69%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5494'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -4.19%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.000008023'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.51%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.45'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.56%  '
$ws.Range("D17").Value = '26.163.54'
$ws.Range("E17").Value = '  -0.24%  '
$ws.Range("E18").Value = '  -0.31%  '
$ws.Range("E19").Value = '  -2.19%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '187.18'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.29%  '
$ws.Range("E21").Value = '  -3.75%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.213'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.36%  '
$ws.Range("E23").Value = '  -0.20%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '150.15'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.89%  '
$ws.Range("E25").Value = '  -1.64%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.494'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.00%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.83'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.26%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.06304'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.59%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.350'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.08%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.284'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.40%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.523'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.12%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.415'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -4.23%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.646'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.93%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.006'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.45%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.6064'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.70%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.404'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.70%  '
$ws.Range("E37").Value = '  +0.41%  '
$ws.Range("D38").Value = '1.113.03'
$ws.Range("E38").Value = '  +1.58%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.115'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.89%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01615'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.05%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8644'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.13%  '
$ws.Range("E42").Value = '  -0.55%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '100.53'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.08%  '
$ws.Range("D44").Value = '1.822.89'
$ws.Range("E44").Value = '  -0.47%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00000000110'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.32%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '55.58'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.10%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.9989'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.59%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.091'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.27%  '
$ws.Range("E49").Value = '  -0.63%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4241'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.87%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '5.931'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.10%  '
